$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: shared-string labels "a".."e" for rows 1-5 ---
$ws.Range("E1").Value = "a"
$ws.Range("E2").Value = "b"
$ws.Range("E3").Value = "c"
$ws.Range("E4").Value = "d"
$ws.Range("E5").Value = "e"

# --- C1: value becomes 1 (was a time-of-day fraction) ---
$ws.Range("C1").Value = 1

# --- C2: value bumped by one day, formatted with the sheetjs.ssf "Long Date" format ---
$ws.Range("C2").Value = 45639
$ws.Range("C2").NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'

# --- Column widths: A becomes an explicit (custom) width, C widens to fit the long date text ---
$ws.Columns("A").ColumnWidth = 10.142857142857142
$ws.Columns("C").ColumnWidth = 23.428571428571427

# --- Selection moves to E6 (just past the newly filled column E) ---
$ws.Range("E6").Select() | Out-Null
